$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.530.52'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.52%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.479.78'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.64%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.17%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.14'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.49%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.56'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.96%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.90%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.22%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.507'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.82%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.72'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.36%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0786'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.83%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.110'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.28%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.860.59'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.67%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.87'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.89%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.25'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +9.58%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.496.70'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.87%  '

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.74%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.529.18'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.69%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.58'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.21%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0944'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.44%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.39'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +5.77%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.24'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.26%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.17'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.53%  '

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.54%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.05%  '

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.37%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.78'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.81%  '

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.13%  '

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.58%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.81'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.78%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '158.03'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +3.94%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.45'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.87%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.58'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.61%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0755'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.42%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.48'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.58%  '

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.56%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.89'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -5.36%  '

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.65%  '

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.04%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.08'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.87%  '

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.30%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.46'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.47%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.961.69'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.29%  '

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.65%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.95'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.16%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.96'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.36%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.717.40'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.46%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '97.69'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.75%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '68.21'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.48%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.37'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.14%  '
